$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct per-cell updates, matching the new scraped cryptos data.
# Numeric-looking Price values need an explicit Text format so Excel keeps them
# as strings (matching the original "inlineStr" cell content) instead of coercing
# them to numbers; the Style reset afterwards clears the temporary formatting so
# the cell keeps its original (unstyled) appearance.

$ws.Range("D2").Value = "20.233.61"
$ws.Range("E2").Value = "  +1.41%  "

$ws.Range("D3").Value = "1.442.61"
$ws.Range("E3").Value = "  +2.51%  "

$ws.Range("E4").Value = "  +0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9183"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.49%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3075"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "38.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.016"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06487"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9984"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.323"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.029"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001008"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").Value = "1.440.91"
$ws.Range("E17").Value = "  +2.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9328"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05622"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.384"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.237"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("D25").Value = "20.261.85"
$ws.Range("E25").Value = "  +1.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.040"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("D29").Value = "1.593.92"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.974"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7888"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.814"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.88%  "

$ws.Range("E34").Value = "  +0.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.461"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05764"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01981"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1843"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9297"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.968"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -16.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5177"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.475"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "115.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5077"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.723"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06381"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9780"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.38%  "

# Rows 37 and 38 swap places: "InternetComputer(DFINITY)" moves down to row 38
# and "TrustWalletToken" moves up to row 37, each with refreshed price/volume data.
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.132"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.01%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.627"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.02%  "
